# Apply cryptos list update (prices/volumes refresh + one row reorder)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number (e.g. "223.71") must be
# forced to Text format first, otherwise Excel would silently convert them
# to numeric values (and normalise things like trailing zeros).
$textCells = @(
    "D5", "D6", "D8", "D9", "D10", "D11", "D13", "D15", "D17", "D18",
    "D19", "D21", "D23", "D24", "D25", "D26", "D27", "D30", "D31", "D32",
    "D34", "D38", "D39", "D40", "D41", "D43", "D44", "D45", "D49", "D50"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "34.547.51"
$ws.Range("E2").Value = "  +1.46%  "
$ws.Range("D3").Value = "1.800.42"
$ws.Range("E3").Value = "  +0.89%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "223.71"
$ws.Range("E5").Value = "  -1.34%  "
$ws.Range("D6").Value = "0.551"
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "32.25"
$ws.Range("E8").Value = "  +3.03%  "
$ws.Range("D9").Value = "0.288"
$ws.Range("E9").Value = "  +2.78%  "
$ws.Range("D10").Value = "0.0706"
$ws.Range("E10").Value = "  +7.36%  "
$ws.Range("D11").Value = "0.0930"
$ws.Range("E11").Value = "  -0.01%  "
$ws.Range("D12").Value = "2.061.44"
$ws.Range("E12").Value = "  +1.08%  "
$ws.Range("D13").Value = "11.04"
$ws.Range("E13").Value = "  -3.11%  "
$ws.Range("D14").Value = "1.803.78"
$ws.Range("E14").Value = "  +1.12%  "
$ws.Range("D15").Value = "0.640"
$ws.Range("E15").Value = "  +0.99%  "
$ws.Range("D16").Value = "34.562.95"
$ws.Range("E16").Value = "  +1.55%  "
$ws.Range("D17").Value = "4.28"
$ws.Range("E17").Value = "  +1.29%  "
$ws.Range("D18").Value = "69.04"
$ws.Range("E18").Value = "  -0.46%  "
$ws.Range("D19").Value = "251.12"
$ws.Range("E19").Value = "  -0.65%  "
$ws.Range("D20").Value = "0.0₃0796"
$ws.Range("E20").Value = "  +7.39%  "
$ws.Range("D21").Value = "11.03"
$ws.Range("E21").Value = "  +5.79%  "
$ws.Range("E22").Value = "  -0.23%  "
$ws.Range("D23").Value = "4.24"
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("D24").Value = "2.16"
$ws.Range("E24").Value = "  +1.09%  "
$ws.Range("D25").Value = "161.52"
$ws.Range("E25").Value = "  +2.84%  "
$ws.Range("D26").Value = "16.35"
$ws.Range("E26").Value = "  -1.24%  "
$ws.Range("D27").Value = "7.12"
$ws.Range("E27").Value = "  +1.35%  "
$ws.Range("E28").Value = "  -0.17%  "
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("D30").Value = "543.81"
$ws.Range("E30").Value = "  +947.99%  "
$ws.Range("D31").Value = "0.0525"
$ws.Range("E31").Value = "  +1.48%  "
$ws.Range("D32").Value = "3.78"
$ws.Range("E32").Value = "  -0.51%  "
$ws.Range("E33").Value = "  -0.27%  "
$ws.Range("D34").Value = "3.60"
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("E35").Value = "  +2.14%  "
$ws.Range("D36").Value = "1.430.19"
$ws.Range("E36").Value = "  -1.34%  "
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("D38").Value = "0.639"
$ws.Range("E38").Value = "  +1.87%  "
$ws.Range("D39").Value = "0.0190"
$ws.Range("E39").Value = "  +1.66%  "
$ws.Range("D40").Value = "84.39"
$ws.Range("E40").Value = "  +1.30%  "
$ws.Range("D41").Value = "0.955"
$ws.Range("E41").Value = "  +6.41%  "
$ws.Range("E42").Value = "  -1.11%  "
$ws.Range("D43").Value = "2.34"
$ws.Range("E43").Value = "  -0.24%  "
$ws.Range("D44").Value = "2.15"
$ws.Range("E44").Value = "  +3.69%  "
$ws.Range("D45").Value = "6.03"
$ws.Range("E45").Value = "  +5.02%  "
$ws.Range("E46").Value = "  -0.82%  "
$ws.Range("E47").Value = "  -2.24%  "
$ws.Range("B48").Value = "RocketPoolETH"
$ws.Range("C48").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D48").Value = "1.955.96"
$ws.Range("E48").Value = "  +0.86%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").Value = "12.26"
$ws.Range("E49").Value = "  +2.80%  "
$ws.Range("D50").Value = "106.18"
$ws.Range("E50").Value = "  +8.69%  "
$ws.Range("E51").Value = "  -0.04%  "
